$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "아디다스 알파바운스슬라이드 실내화 슬리퍼 BA8775 B41720 FZ0387 FZ0388"
$ws.Range("B4").Value = "아디다스 듀라모 슬라이드 슬리퍼 FY8786"
$ws.Range("B6").Value = "크록스 공용 바야밴드 슬라이드 슬리퍼 20SUSL205733"
$ws.Range("B7").Value = "아디다스 남여공용 슬리퍼 알파바운스 슬라이드 ALPHABOUNCE SLIDE FZ0387"
$ws.Range("B8").Value = "SNRD 군납용슬리퍼 보급 군용 군대 군인 군용슬리퍼"
$ws.Range("B9").Value = "아디다스 아딜렛 CF 로고 슬리퍼 CG3425"
$ws.Range("B10").Value = "크록스 바야밴드 슬라이드 슬리퍼 205392-066"
$ws.Range("B12").Value = "푸마 슬리퍼 스커프 패딩 방한 384945-01"
$ws.Range("B13").Value = "크록스 슬라이드 사무실 슬리퍼 205733"
$ws.Range("B14").Value = "크록스 크록밴드Ⅲ 205733 슬리퍼 3종 모음 택1"
$ws.Range("B15").Value = "뉴발란스 NEW BALANCE 슬리퍼 unisex SD1501"
$ws.Range("B16").Value = "언더아머 UA 남성 플레이메이커 픽스드 스트랩 슬리퍼 3000061"
$ws.Range("B17").Value = "아디다스 슬리퍼 듀라모 알파바운스 BA8775 B41720"
$ws.Range("B18").Value = "푸마 로얄캣 컴포트 벨크로 슬리퍼 372280-02"
$ws.Range("B19").Value = "콜맨 콤테 방한 슬리퍼 남성 털 실내화"
$ws.Range("B20").Value = "나이키 베나시 JDI 슬리퍼 343880-090"
$ws.Range("B21").Value = "푸마 공용 스커프 털 슬리퍼 Scuff Sherpa 셰르파 슬리퍼 384943"
$ws.Range("B23").Value = "나이키 카와 샤워 슬리퍼 832528-001"
$ws.Range("B24").Value = "아디다스 슬리퍼 아딜렛 컴포트 EG1850"
$ws.Range("B25").Value = "소프트 남성 남자 사무실 슬리퍼 실내화 샌들"
$ws.Range("B26").Value = "남성 여성 커플 겨울 털 슬리퍼 실내화 WT9771"
$ws.Range("B28").Value = "노르딕 겨울 캠핑 슈즈 방한 남자 털 슬리퍼"
$ws.Range("B29").Value = "휠라 TAPER 테이퍼 슬라이드 슬리퍼 1SM00559"
$ws.Range("B31").Value = "디스커버리익스페디션 남여공용 레스터 뮬 VE 겨울 슬리퍼 DXSH4A061"
$ws.Range("B32").Value = "나이키 슬리퍼 오프코트 슬라이드 BQ4639-001"
$ws.Range("B33").Value = "아디다스 남여공용 슬리퍼 아딜렛컴포트 FX4293"
$ws.Range("B35").Value = "호주 어그 UGG 블레드 남자 양털슬리퍼"
$ws.Range("B36").Value = "뉴발란스 남성 슬리퍼 SD1101HWB"
$ws.Range("B37").Value = "슈펜 전황일 콜라보 남성 캠핑 리커버리 샌들 슬리퍼 AFDU79S07"
$ws.Range("B39").Value = "휠라 드리프터 테이피테잎 슬리퍼 1SM00561"
$ws.Range("B42").Value = "언더아머 이그나이트 5 슬리퍼 1287318-100"
$ws.Range("B43").Value = "아디다스 아딜렛 아쿠아 슬리퍼 F35550"
$ws.Range("B44").Value = "버켄스탁 보스턴 타우페 SFB 남성 슬리퍼 REGULAR 560771"
$ws.Range("B45").Value = "페이퍼플레인 남녀공용 방한화 털슬리퍼 실내화 PP1508"
$ws.Range("B46").Value = "크록스 슬라이드 크록밴드 3 슬리퍼 205733-462"
$ws.Range("B47").Value = "언더아머 슬리퍼 이그나이트 3022711"
$ws.Range("B48").Value = "푸마 슬리퍼 Divecat V2 369400"
$ws.Range("B49").Value = "어그 남여공용 양모 슬리퍼 2색상"
$ws.Range("B50").Value = "아디다스 슬리퍼 A23- 듀라모 슬라이드 FY6034"
$ws.Range("B53").Value = "베어파우 ROGER 양털슬리퍼 mens"
$ws.Range("B55").Value = "나이키 슬리퍼 카와 슬라이드 832646-010 남성"
$ws.Range("B56").Value = "나이키 빅토리원 슬라이드 슬리퍼 남자 여자 학생 실내화 CN9675-100"
$ws.Range("B57").Value = "수부 수부 방한 패딩 슬리퍼 다운 샌들 SUBU-BK"
$ws.Range("B58").Value = "뉴발란스 1501 슬라이드 벨크로 슬리퍼 남자 여자 학교 찍찍이 SD1501ILG"
$ws.Range("B59").Value = "아디다스 아딜렛 아쿠아 슬라이드 슬리퍼 F35543"
$ws.Range("B61").Value = "국내산 버켄 보스턴 털슬리퍼 겨울 블로퍼 실내화 사무실용 남녀공용 커플"
$ws.Range("B62").Value = "코스트코 털크록스 크록스털신 가을겨울 발편한 따뜻한 주방 작업 사무실 털슬리퍼 280"
$ws.Range("B63").Value = "아디다스 아딜렛 클로그 슬리퍼 신발 FY8970"
$ws.Range("B65").Value = "오즈웨어 메이슨 남성 슬리퍼 OB687"
$ws.Range("B66").Value = "호카오네오네 호카 남성 오라 리커버리 슬리퍼 BBLC 1099674"
$ws.Range("B67").Value = "내셔널지오그래픽 신발 다이브 플러스 V2 슬리퍼 N211AFW150"
$ws.Range("B69").Value = "나이키 250 오프코트 슬라이드 슬리퍼 BQ4639-012"
$ws.Range("B70").Value = "아디다스 알파바운스 슬라이드 슬리퍼모음 외3종 BA8775"
$ws.Range("B71").Value = "뉴발란스 NB 슬라이드 슬리퍼 남자 여자 실내화 신발 SD1101HBB"
$ws.Range("B72").Value = "언더아머 UA 코어 PTH 슬리퍼 3021286"
$ws.Range("B73").Value = "휠라 드리프터 슬라이드 슬리퍼 WWT WNV 1SM00560"
$ws.Range("B74").Value = "우포스 슬리퍼 OOAHH BLACK"
$ws.Range("B78").Value = "아디다스 아딜렛 컴포트 슬리퍼 CG3425"
$ws.Range("B79").Value = "버켄스탁 취리히 타우페 레귤러 슬리퍼 1009532"
$ws.Range("B81").Value = "크록스 2종균일가 크록밴드 II 슬라이드 슬리퍼 204108"
$ws.Range("B82").Value = "내셔널지오그래픽 다이브 플러스 빅로고 슬리퍼 N205AFW900"
$ws.Range("B83").Value = "나이키 슬리퍼 빅토리 원 슬라이드 흰파 CN9675-102"
$ws.Range("B84").Value = "나이키 슬리퍼 슬라이드 빅토리 원 CZ5478-001"
$ws.Range("B88").Value = "뉴발란스 NC02 슬리퍼 SD1101FZK"
$ws.Range("B89").Value = "아디다스 남녀공용 슬리퍼 듀라모 에스엘 슬라이드 DURAMO SLIDE FY8786"
$ws.Range("B90").Value = "뉴발란스 토앤토 플립플랍 아이보리 리커버리 쪼리 슬리퍼 SD5601GIV"
$ws.Range("B91").Value = "아디다스 아딜렛 CF 슬리퍼 블랙 CG3427"
$ws.Range("B93").Value = "나이키 슬리퍼 베네시 JDI 맨즈 343880-090"
$ws.Range("B94").Value = "남녀공용 기모안감 털슬리퍼 6컬러 DS-AL641"
$ws.Range("B95").Value = "아디다스 슬리퍼 NQB 알파바운스 슬라이드 슬리퍼 B41720"
$ws.Range("B96").Value = "나이키 슬리퍼 A4- 빅토리 원 샤워 CZ5478-001"
$ws.Range("B97").Value = "뉴발란스 슬리퍼 1501 남여공용 쿠셔닝 SD1501IWT"
$ws.Range("B98").Value = "MLB LA 다저스BS 방한 슬리퍼 NW"
$ws.Range("B99").Value = "클럽 남성 남자 사무실 슬리퍼 실내화 샌들 비치"
$ws.Range("B100").Value = "노스페이스 뮬 슬리퍼 플리스 A NS93M60A"
$ws.Range("B101").Value = "MLB 방한 EVA슬리퍼"
